{"js": "// Commit: \"bea1 c on dos lineas\"\n//\n// The original document has a single paragraph containing the text\n// \"rererer\". The edit:\n//   1. Splits \"rererer\" into two runs \"R\" + \"ererer\" (i.e. capitalises the\n//      first letter while keeping it as two separate <w:r> runs) and wraps\n//      it with a spell-check proofing-error marker pair\n//      (<w:proofErr w:type=\"spellStart\"/> ... <w:proofErr w:type=\"spellEnd\"/>).\n//   2. Adds a brand-new second paragraph with the text \"Beatrizv fg\", built\n//      out of runs \"B\" + \"eatriz\" + \"v\" (wrapped in its own proofErr pair)\n//      followed by a run containing \" fg\". The _GoBack bookmark that used to\n//      sit at the end of paragraph 1 now sits at the end of this new\n//      paragraph 2 (i.e. at the end of the document, same as before).\n//\n// We reproduce this precisely (including the literal <w:proofErr/> markers,\n// which Word's spell checker would normally stamp in and which the Office.js\n// object model otherwise has no dedicated property for) by replacing the\n// whole first paragraph's range with the equivalent literal OOXML, wrapped\n// in the \"flat OPC\" package format required by Range.insertOoxml.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\nconst documentXml =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n      '<w:p>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>R</w:t></w:r>' +\n        '<w:r><w:t>ererer</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n      '</w:p>' +\n      '<w:p>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>B</w:t></w:r>' +\n        '<w:r><w:t>eatriz</w:t></w:r>' +\n        '<w:r><w:t>v</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> fg</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n      '</w:p>' +\n    '</w:body>' +\n  '</w:document>';\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n      '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n          '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' + documentXml + '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\n// The paragraph's whole range (text + its end-of-paragraph mark, which is\n// also where the _GoBack bookmark lives) gets replaced by the two\n// paragraphs above in one shot.\nconst wholeRange = firstParagraph.getRange(\"Whole\");\nwholeRange.insertOoxml(flatOpcXml, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Commit: \"bea1 c on dos lineas\"\n#\n# The original document has a single paragraph containing the text\n# \"rererer\". The edit:\n#   1. Splits \"rererer\" into two runs \"R\" + \"ererer\" (i.e. capitalises the\n#      first letter while keeping it as two separate runs) and wraps it with\n#      a spell-check proofing-error marker pair (proofErr spellStart/spellEnd).\n#   2. Adds a brand-new second paragraph with the text \"Beatrizv fg\", built\n#      out of runs \"B\" + \"eatriz\" + \"v\" (wrapped in its own proofErr pair)\n#      followed by a run containing \" fg\". The _GoBack bookmark that used to\n#      sit at the end of paragraph 1 now sits at the end of this new\n#      paragraph 2 (i.e. still at the end of the document).\n#\n# We reproduce this precisely (including the literal <w:proofErr/> markers)\n# by replacing the whole first paragraph's Range with the equivalent literal\n# OOXML via Range.InsertXML, using the \"flat OPC\" package format.\n\n$d = $word.ActiveDocument\n\n$documentXml = '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>R</w:t></w:r><w:r><w:t>ererer</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>B</w:t></w:r><w:r><w:t>eatriz</w:t></w:r><w:r><w:t>v</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> fg</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p></w:body></w:document>'\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + $documentXml + '</pkg:xmlData></pkg:part></pkg:package>'\n\n# Paragraph 1's whole Range (its text plus its end-of-paragraph mark, where\n# the _GoBack bookmark lives) gets replaced by the two paragraphs above.\n$r = $d.Paragraphs(1).Range\n$r.InsertXML($flatOpcXml) | Out-Null\n"}
